$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A57").NumberFormat = "@"
$ws.Range("A57").Value = "10/28/2025"
$ws.Range("B57").Value = 0.1865144124783668
$ws.Range("C57").Value = 0.8134855875216332
